$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1
$ws.Range("F1").Value = "sd_ETR"

# Add sd_ETR values for each data row (rows 2-17)
$values = @(
    0.0254991886750858,
    13.0484203546695,
    19.3580536394718,
    3.42081448184145,
    0,
    25.4418999660482,
    21.4200597475962,
    7.34017452807367,
    0.0324652687876046,
    59.2576446341752,
    44.713832935325,
    3.571556864073,
    0.0911456381341045,
    53.8508294947713,
    42.4242775391641,
    0
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 6).Value = $v
    $row++
}
